$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "Role" (column D), shifting Role/Operation
# one column to the right (D->E, E->F). Excel copies the left neighbour's
# formatting into the new column as part of the insert.
$ws.Columns("D:D").Insert()

# New column header
$ws.Range("D1").Value = "User License"

# Match the new column's width to the (hyperlink) Email column next to it
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# Row 4's inserted cell sits in the thick-bordered "box" that used to wrap
# only column C; splitting it removes the border on the new interior cell
# while keeping the fill/alignment.
$ws.Range("D4").Style = "Hyperlink"
$ws.Range("D4").Interior.Color = $ws.Range("C4").Interior.Color
$ws.Range("D4").VerticalAlignment = -4108
$ws.Range("D4").WrapText = $true
$ws.Range("D4").Borders.LineStyle = 0

# Match the selection left behind after typing the new header
$ws.Range("D1").Select()
